$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -5.028042793650829
$ws.Range("C3").Value = -1.189874591947448
$ws.Range("C4").Value = -0.0615282423017409
$ws.Range("C5").Value = -0.4005760218502046
$ws.Range("C6").Value = 0.01500614343439477
$ws.Range("C7").Value = 0.09915277476007525
$ws.Range("C8").Value = 0.1276481233015081
$ws.Range("C9").Value = 0.02713608500475133
$ws.Range("C10").Value = 0.03198635883432693
$ws.Range("C11").Value = 0.00527679405225358
$ws.Range("C12").Value = 0.03738325939343457
